{"js": "const replacements = [\n  [\"Dna. Ioana Mirea\", \"Dl. Adrian Terec\"],\n  [\"\u015eef Birou Aprovizionare\", \"Director Executiv\"],\n  [\"Birou Achizi\u0163ii\", \"Management\"],\n  [\"S.C. Automatica S.A.\", \"S.C. EnergoBIT S.R.L.\"],\n  [\"Bdul. Voluntari Nr. 108 Bis, Voluntari\", \"Strada T\u0103ietura Turcului Nr. 47, Cluj-Napoca\"],\n  [\"Tel: +40 (372) 058 100*107\", \"Tel: +40 (264) 207 544\"],\n  [\"Fax: +40 (372) 058 101\", \"Fax: +40 (264) 207 555\"],\n  [\"Mobil: +40 (729) 035 164\", \"Mobil: +40 (723) 658 773\"],\n  [\"ioana.mirea@automatica.ro\", \"adrian.terec@energobit.com\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Dna. Ioana Mirea\"; Replace = \"Dl. Adrian Terec\" },\n    @{ Find = \"\u015eef Birou Aprovizionare\"; Replace = \"Director Executiv\" },\n    @{ Find = \"Birou Achizi\u0163ii\"; Replace = \"Management\" },\n    @{ Find = \"S.C. Automatica S.A.\"; Replace = \"S.C. EnergoBIT S.R.L.\" },\n    @{ Find = \"Bdul. Voluntari Nr. 108 Bis, Voluntari\"; Replace = \"Strada T\u0103ietura Turcului Nr. 47, Cluj-Napoca\" },\n    @{ Find = \"Tel: +40 (372) 058 100*107\"; Replace = \"Tel: +40 (264) 207 544\" },\n    @{ Find = \"Fax: +40 (372) 058 101\"; Replace = \"Fax: +40 (264) 207 555\" },\n    @{ Find = \"Mobil: +40 (729) 035 164\"; Replace = \"Mobil: +40 (723) 658 773\" },\n    @{ Find = \"ioana.mirea@automatica.ro\"; Replace = \"adrian.terec@energobit.com\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
